$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.484.63"
$ws.Range("E2").Value = "  +0.47%  "

$ws.Range("D3").Value = "2.488.10"
$ws.Range("E3").Value = "  -2.29%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.06%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "588.34"
$ws.Range("E5").Value = "  -0.59%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "171.71"
$ws.Range("E6").Value = "  -1.33%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.523"
$ws.Range("E8").Value = "  -1.11%  "

$ws.Range("D9").Value = "2.484.22"
$ws.Range("E9").Value = "  -2.47%  "

$ws.Range("E10").Value = "  -0.73%  "

$ws.Range("E11").Value = "  +1.15%  "

$ws.Range("E12").Value = "  -1.08%  "

$ws.Range("E13").Value = "  -2.23%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.17"
$ws.Range("E14").Value = "  -2.60%  "

$ws.Range("D15").Value = "2.944.63"
$ws.Range("E15").Value = "  -2.22%  "

$ws.Range("E16").Value = "  -0.71%  "

$ws.Range("D17").Value = "67.288.54"
$ws.Range("E17").Value = "  +0.46%  "

$ws.Range("D18").Value = "2.495.48"
$ws.Range("E18").Value = "  -2.40%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.63"
$ws.Range("E19").Value = "  +2.54%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.78"
$ws.Range("E20").Value = "  -3.54%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "366.16"
$ws.Range("E21").Value = "  +3.08%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.13"
$ws.Range("E22").Value = "  -1.59%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.54"
$ws.Range("E23").Value = "  -2.01%  "

$ws.Range("E24").Value = "  -0.01%  "

$ws.Range("E25").Value = "  +1.85%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.91"
$ws.Range("E26").Value = "  -4.11%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.93"
$ws.Range("E27").Value = "  -1.68%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.994"
$ws.Range("E28").Value = "  -0.39%  "

$ws.Range("E29").Value = "  -2.89%  "

$ws.Range("D30").Value = "0.0₃0956"
$ws.Range("E30").Value = "  -3.46%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.25"
$ws.Range("E31").Value = "  +0.60%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "528.25"
$ws.Range("E32").Value = "  -1.68%  "

$ws.Range("E33").Value = "  -3.28%  "

$ws.Range("E34").Value = "  -0.02%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.998"
$ws.Range("E35").Value = "  -0.17%  "

$ws.Range("E36").Value = "  -3.47%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "159.32"
$ws.Range("E37").Value = "  +1.54%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.41"
$ws.Range("E38").Value = "  -4.19%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.72"
$ws.Range("E39").Value = "  +0.20%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.62"
$ws.Range("E40").Value = "  +0.85%  "

$ws.Range("E41").Value = "  -2.48%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.76"
$ws.Range("E42").Value = "  -1.84%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.08"
$ws.Range("E43").Value = "  -1.44%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.00"
$ws.Range("E44").Value = "  +0.13%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.44"
$ws.Range("E45").Value = "  -3.19%  "

$ws.Range("E46").Value = "  -0.82%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "144.71"
$ws.Range("E47").Value = "  -3.33%  "

$ws.Range("E48").Value = "  -1.00%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.545"
$ws.Range("E49").Value = "  -2.94%  "

$ws.Range("E50").Value = "  -1.39%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0744"
$ws.Range("E51").Value = "  -2.32%  "
